# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values replacing the old Strike# derived values in column G,
# for data rows 2 through 38 (row 1 is the header row).
$newK = @{
    2  = 4
    3  = 3
    4  = 6
    5  = 4
    6  = 7
    7  = 4
    8  = 6
    9  = 3
    10 = 0
    11 = 8
    12 = 6
    13 = 6
    14 = 3
    15 = 3
    16 = 4
    17 = 4
    18 = 5
    19 = 4
    20 = 3
    21 = 2
    22 = 4
    23 = 5
    24 = 2
    25 = 5
    26 = 3
    27 = 5
    28 = 1
    29 = 4
    30 = 4
    31 = 6
    32 = 4
    33 = 2
    34 = 8
    35 = 2
    36 = 4
    37 = 3
    38 = 1
}

foreach ($row in $newK.Keys | Sort-Object) {
    $ws.Range("G$row").Value = $newK[$row]
}
